# Adjust row heights on the active sheet to match the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row gets a bit taller
$ws.Rows.Item(1).RowHeight = 18.75

# Data rows 2 through 39 get taller as well
$ws.Range("A2:A39").EntireRow.RowHeight = 19.5

# Row 77 (near the bottom) shrinks back down
$ws.Rows.Item(77).RowHeight = 17.25
